$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value2 = '23.519.45'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value2 = '  -0.76%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value2 = '1.639.80'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value2 = '  -0.84%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value2 = '0.9981'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value2 = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '1.000'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value2 = '  +0.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value2 = '304.24'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value2 = '  -0.68%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value2 = '0.3786'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value2 = '  +0.26%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value2 = '51.80'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value2 = '  -1.43%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '0.3641'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value2 = '  -0.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value2 = '0.08188'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value2 = '  +0.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value2 = '1.231'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value2 = '  -3.17%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value2 = '0.9982'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value2 = '  -0.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value2 = '22.50'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value2 = '  -2.78%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '6.473'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value2 = '  -3.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value2 = '7.388'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value2 = '  -0.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value2 = '0.00001242'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value2 = '  -2.70%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value2 = '1.636.67'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value2 = '  -1.69%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value2 = '95.14'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value2 = '  -0.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '0.06942'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value2 = '  +0.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '6.596'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value2 = '  -0.26%  '
$ws.Range('E21').Value2 = '  -5.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '0.9992'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value2 = '  +0.07%  '
$ws.Range('E23').Value2 = '  -3.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value2 = '23.512.98'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value2 = '  -0.83%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value2 = '2.511'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value2 = '  +3.90%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value2 = '3.076'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value2 = '  -3.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value2 = '21.14'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value2 = '  -1.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value2 = '151.79'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value2 = '  +0.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value2 = '5.269'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value2 = '  -1.20%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value2 = '133.45'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value2 = '  -3.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value2 = '1.817.64'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value2 = '  -1.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '6.638'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value2 = '  -4.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '2.163'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value2 = '  -7.75%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value2 = '1.051'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value2 = '  +7.86%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value2 = '11.39'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value2 = '  +3.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '0.02767'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value2 = '  -3.62%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value2 = '0.2493'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value2 = '  -3.74%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value2 = '0.08779'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value2 = '  -1.38%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value2 = '0.07138'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value2 = '  -3.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value2 = '6.037'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value2 = '  -5.74%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value2 = '0.7029'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value2 = '  -3.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value2 = '1.343'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value2 = '  -2.68%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '12.20'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value2 = '  -4.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value2 = '15.85'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value2 = '  -4.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '0.6520'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value2 = '  -2.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value2 = '0.9996'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value2 = '  +0.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value2 = '2.283'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value2 = '  -4.09%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '3.968'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value2 = '  -1.49%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value2 = '0.07978'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value2 = '  -0.88%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '127.23'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value2 = '  -1.07%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value2 = '1.196'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value2 = '  -2.73%  '
